# Mark additional "Greedy" (and a couple of BST) rows as done in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yesRows = @(216, 218, 243, 245, 246, 247, 250, 251, 252, 253, 255, 256, 257, 260, 269)
foreach ($r in $yesRows) {
    $ws.Range("C$r").Value = "yes"
}

# Row 244 was marked with a capitalized "Yes" instead.
$ws.Range("C244").Value = "Yes"

# Update the saved view state (scroll position / active selection) to match
# where the author ended up working.
$ws.Application.Goto($ws.Range("B265"), $true)
$window = $excel.ActiveWindow
$window.ScrollRow = 233
$window.ScrollColumn = 1
$ws.Range("B265").Select()
